$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.922.43'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '1.551.22'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.69'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.09'
$ws.Range("E8").Value = '  +3.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.246'
$ws.Range("E9").Value = '  -0.20%  '
$ws.Range("E10").Value = '  +0.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0857'
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = '1.772.92'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").Value = '1.560.19'
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("E15").Value = '  +1.51%  '
$ws.Range("D16").Value = '26.931.07'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.67'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.42'
$ws.Range("E18").Value = '  +1.89%  '
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.27'
$ws.Range("E20").Value = '  +1.20%  '
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.05'
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.21'
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.15'
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.61'
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.90'
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("E30").Value = '  +2.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.09'
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("D33").Value = '1.422.34'
$ws.Range("E33").Value = '  +4.41%  '
$ws.Range("E34").Value = '  +4.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.972'
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0165'
$ws.Range("E38").Value = '  +0.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.520'
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.807'
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.75'
$ws.Range("E41").Value = '  +5.36%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  +4.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.993'
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.31'
$ws.Range("E45").Value = '  +1.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.76'
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("D47").Value = '1.686.98'
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.68'
$ws.Range("E48").Value = '  +1.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0519'
$ws.Range("E49").Value = '  +2.57%  '
$ws.Range("E50").Value = '  +3.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0951'
$ws.Range("E51").Value = '  +0.27%  '
